$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (D last_edited_time, AC properties.Đầy đủ.number, AF properties.Tổng công.number)
# '-' entries mean that column is unchanged for that row.

$ws.Range("D2").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC2").Value = 17
$ws.Range("AF2").Value = 17

$ws.Range("D3").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC3").Value = 17
$ws.Range("AF3").Value = 17

$ws.Range("D4").Value = "2024-07-18T15:58:00.000Z"

$ws.Range("D5").Value = "2024-07-18T15:58:00.000Z"

$ws.Range("D6").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC6").Value = 17
$ws.Range("AF6").Value = 17

$ws.Range("D7").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC7").Value = 17
$ws.Range("AF7").Value = 17

$ws.Range("D8").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC8").Value = 17
$ws.Range("AF8").Value = 18

$ws.Range("D9").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC9").Value = 16
$ws.Range("AF9").Value = 17

$ws.Range("D10").Value = "2024-07-18T15:58:00.000Z"

$ws.Range("D11").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC11").Value = 18
$ws.Range("AF11").Value = 19.5

$ws.Range("D12").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC12").Value = 15
$ws.Range("AF12").Value = 17

$ws.Range("D13").Value = "2024-07-18T15:58:00.000Z"

$ws.Range("D14").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC14").Value = 17
$ws.Range("AF14").Value = 17.5

$ws.Range("D15").Value = "2024-07-18T15:58:00.000Z"

$ws.Range("D16").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC16").Value = 11
$ws.Range("AF16").Value = 16

$ws.Range("D17").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC17").Value = 17
$ws.Range("AF17").Value = 17.5

$ws.Range("D18").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC18").Value = 17
$ws.Range("AF18").Value = 17

$ws.Range("D19").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC19").Value = 17
$ws.Range("AF19").Value = 18

$ws.Range("D20").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("AC20").Value = 18
$ws.Range("AF20").Value = 18
